# Auto-generated edit script
# Applies updated profit-calculation values to the Odin_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 313.66666
$ws.Range("I41").Value = 313.66666
$ws.Range("K41").Value = 313.66666
$ws.Range("M41").Value = 126.33334

$ws.Range("H42").Value = 562.3333
$ws.Range("I42").Value = 943.6667
$ws.Range("K42").Value = 2831.0001
$ws.Range("M42").Value = -2601.0001

$ws.Range("H70").Value = 2099.4092
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 2232.611
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 6697.833
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -7237.833

$ws.Range("H73").Value = 2099.4092
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 2232.611
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 6697.833
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -8569.832999999999

$ws.Range("H76").Value = 200003170
$ws.Range("I76").Value = 250002960
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 250002960
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -250002645
$ws.Range("N76").Value = -4630

$ws.Range("H79").Value = 200003170
$ws.Range("I79").Value = 250002960
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 250002960
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -250001868
$ws.Range("N79").Value = -6184

$ws.Range("H115").Value = 1180.125
$ws.Range("I115").Value = 1205.8572
$ws.Range("K115").Value = 3617.5716
$ws.Range("M115").Value = -2050.5716

$ws.Range("H132").Value = 465080.2
$ws.Range("I132").Value = 539677
$ws.Range("K132").Value = 1619031
$ws.Range("M132").Value = -1616501

$ws.Range("H135").Value = 4982.793
$ws.Range("I135").Value = 1330.3
$ws.Range("K135").Value = 11972.7
$ws.Range("M135").Value = -9437.699999999999

$ws.Range("H138").Value = 2939.365
$ws.Range("J138").Value = 4973.242
$ws.Range("L138").Value = 14919.726
$ws.Range("N138").Value = -25199.726

$ws.Range("H141").Value = 2526.353
$ws.Range("I141").Value = 2421.0908
$ws.Range("K141").Value = 7263.2724
$ws.Range("M141").Value = -2083.2724


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1340889
$ws.Range("I32").Value = 874.4769
$ws.Range("K32").Value = 874.4769
$ws.Range("M32").Value = -587.4769

$ws.Range("H63").Value = 4986.4287
$ws.Range("I63").Value = 902.5
$ws.Range("J63").Value = 6620
$ws.Range("K63").Value = 902.5
$ws.Range("L63").Value = 6620
$ws.Range("M63").Value = -216.5
$ws.Range("N63").Value = -7992

$ws.Range("H66").Value = 4986.4287
$ws.Range("I66").Value = 902.5
$ws.Range("J66").Value = 6620
$ws.Range("K66").Value = 4512.5
$ws.Range("L66").Value = 33100
$ws.Range("M66").Value = -1080.5
$ws.Range("N66").Value = -39964

$ws.Range("H132").Value = 849365.4
$ws.Range("I132").Value = 948421.3
$ws.Range("K132").Value = 2845263.9
$ws.Range("M132").Value = -2842733.9


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 49994.5
$ws.Range("J44").Value = 49994.5
$ws.Range("L44").Value = 49994.5
$ws.Range("N44").Value = -50988.5

$ws.Range("H80").Value = 14504706
$ws.Range("I80").Value = 1407.1
$ws.Range("J80").Value = 25661090
$ws.Range("K80").Value = 1407.1
$ws.Range("L80").Value = 25661090
$ws.Range("M80").Value = -409.0999999999999
$ws.Range("N80").Value = -25663086

$ws.Range("H83").Value = 14504706
$ws.Range("I83").Value = 1407.1
$ws.Range("J83").Value = 25661090
$ws.Range("K83").Value = 7035.5
$ws.Range("L83").Value = 128305450
$ws.Range("M83").Value = -2043.5
$ws.Range("N83").Value = -128315434

$ws.Range("H134").Value = 1771663.8
$ws.Range("I134").Value = 2652954
$ws.Range("K134").Value = 7958862
$ws.Range("M134").Value = -7956327


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1090606
$ws.Range("I22").Value = 1716403.4
$ws.Range("J22").Value = 107210.14
$ws.Range("K22").Value = 1716403.4
$ws.Range("L22").Value = 107210.14
$ws.Range("M22").Value = -1716053.4
$ws.Range("N22").Value = -107910.14

$ws.Range("H45").Value = 63744.5
$ws.Range("J45").Value = 63989.5
$ws.Range("L45").Value = 63989.5
$ws.Range("N45").Value = -65175.5

$ws.Range("H86").Value = 8871.352999999999
$ws.Range("I86").Value = 4753.909
$ws.Range("K86").Value = 4753.909
$ws.Range("M86").Value = -3630.909

$ws.Range("H89").Value = 8871.352999999999
$ws.Range("I89").Value = 4753.909
$ws.Range("K89").Value = 23769.545
$ws.Range("M89").Value = -18153.545

$ws.Range("H132").Value = 7388.2915
$ws.Range("I132").Value = 5801.1055
$ws.Range("K132").Value = 17403.3165
$ws.Range("M132").Value = -14873.3165


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 11782.667
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 11782.667
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 35348.001
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -35978.001

$ws.Range("H73").Value = 11782.667
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 11782.667
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 35348.001
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -37532.001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4822.963
$ws.Range("I102").Value = 4074.4285
$ws.Range("K102").Value = 4074.4285
$ws.Range("M102").Value = -2452.4285


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2515
$ws.Range("I68").Value = 2354.7
$ws.Range("J68").Value = 3049.3333
$ws.Range("K68").Value = 2354.7
$ws.Range("L68").Value = 3049.3333
$ws.Range("M68").Value = -1605.7
$ws.Range("N68").Value = -4547.3333

$ws.Range("H71").Value = 2515
$ws.Range("I71").Value = 2354.7
$ws.Range("J71").Value = 3049.3333
$ws.Range("K71").Value = 11773.5
$ws.Range("L71").Value = 15246.6665
$ws.Range("M71").Value = -8029.5
$ws.Range("N71").Value = -22734.6665

$ws.Range("H132").Value = 4508.4165
$ws.Range("I132").Value = 4665.6
$ws.Range("K132").Value = 13996.8
$ws.Range("M132").Value = -11466.8


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1746.8572
$ws.Range("I100").Value = 879.4
$ws.Range("J100").Value = 2228.7778
$ws.Range("K100").Value = 1758.8
$ws.Range("L100").Value = 4457.5556
$ws.Range("M100").Value = -1217.8
$ws.Range("N100").Value = -5539.5556

